# Updates the cryptos list (Sheet1, rows 2-51) per the Sep 21 2024 refresh:
# new Price/Volume(1h) figures, plus four coins re-ranked (rows 25-28, 47-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Columns D (Price) and E (Volume) store plain text (e.g. "63.125.43",
    # "0.585", "  +0.24%  "). A leading apostrophe forces Excel to keep the
    # assignment as literal text instead of parsing it into a number; the
    # Style reset afterwards drops the resulting quote-prefix formatting so
    # the cell keeps its original (default) style.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "63.168.20"
Set-TextCell "E2" "  +0.64%  "
Set-TextCell "D3" "2.572.52"
Set-TextCell "E3" "  +4.78%  "
Set-TextCell "D5" "571.21"
Set-TextCell "E5" "  +1.68%  "
Set-TextCell "D6" "147.10"
Set-TextCell "E6" "  +3.47%  "
Set-TextCell "E7" "  -0.06%  "
Set-TextCell "D8" "0.583"
Set-TextCell "E8" "  -0.16%  "
Set-TextCell "D9" "2.574.61"
Set-TextCell "E9" "  +4.92%  "
Set-TextCell "E10" "  +0.78%  "
Set-TextCell "E11" "  -1.22%  "
Set-TextCell "E12" "  +0.07%  "
Set-TextCell "E13" "  +0.66%  "
Set-TextCell "D14" "27.63"
Set-TextCell "E14" "  +3.15%  "
Set-TextCell "D15" "3.033.97"
Set-TextCell "E15" "  +4.81%  "
Set-TextCell "D16" "63.104.69"
Set-TextCell "E16" "  +0.67%  "
Set-TextCell "D17" "0.0000144"
Set-TextCell "E17" "  +2.55%  "
Set-TextCell "D18" "2.578.87"
Set-TextCell "E18" "  +4.94%  "
Set-TextCell "D19" "11.45"
Set-TextCell "E19" "  +2.43%  "
Set-TextCell "D20" "335.41"
Set-TextCell "E20" "  -0.83%  "
Set-TextCell "D21" "4.34"
Set-TextCell "E21" "  +2.60%  "
Set-TextCell "D22" "6.85"
Set-TextCell "E22" "  +1.42%  "
Set-TextCell "E23" "  +0.06%  "
Set-TextCell "D24" "65.34"
Set-TextCell "E24" "  +0.07%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D25" "0.171"
Set-TextCell "E25" "  +0.43%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D26" "1.62"
Set-TextCell "E26" "  +8.91%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D27" "1.00"
Set-TextCell "E27" "  +0.05%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D28" "8.46"
Set-TextCell "E28" "  +5.66%  "
Set-TextCell "D29" "1.49"
Set-TextCell "E29" "  +7.22%  "
Set-TextCell "D30" "7.35"
Set-TextCell "E30" "  +8.77%  "
Set-TextCell "D31" "0.0₃0826"
Set-TextCell "E31" "  +5.03%  "
Set-TextCell "D32" "1.86"
Set-TextCell "E32" "  +1.69%  "
Set-TextCell "D33" "175.56"
Set-TextCell "E33" "  +0.08%  "
Set-TextCell "D34" "1.56"
Set-TextCell "E34" "  +3.87%  "
Set-TextCell "D35" "406.77"
Set-TextCell "E35" "  +10.02%  "
Set-TextCell "D36" "0.401"
Set-TextCell "E36" "  +1.05%  "
Set-TextCell "D37" "19.02"
Set-TextCell "E37" "  +1.76%  "
Set-TextCell "E38" "  -0.01%  "
Set-TextCell "E39" "  +1.83%  "
Set-TextCell "D40" "1.76"
Set-TextCell "E40" "  +4.90%  "
Set-TextCell "E41" "  +0.01%  "
Set-TextCell "D42" "39.37"
Set-TextCell "E42" "  -2.22%  "
Set-TextCell "D43" "153.54"
Set-TextCell "E43" "  +2.97%  "
Set-TextCell "D44" "3.77"
Set-TextCell "E44" "  +2.50%  "
Set-TextCell "D45" "21.04"
Set-TextCell "E45" "  +3.47%  "
Set-TextCell "D46" "0.609"
Set-TextCell "E46" "  +1.90%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D47" "0.0242"
Set-TextCell "E47" "  +7.22%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D48" "0.0528"
Set-TextCell "E48" "  +2.82%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D49" "0.0963"
Set-TextCell "E49" "  +0.75%  "
Set-TextCell "D50" "18.60"
Set-TextCell "E50" "  +4.34%  "
Set-TextCell "E51" "  +1.20%  "
